# "Add files via upload" — the author replaced the old single-column
# consignment-number lookup sheet with a small "TestData" table (and
# renamed the sheet accordingly). Reproduce that on sheet 1, leaving the
# other three sheets (consignmentsummary / precondetails / eventdetails)
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the first sheet: "consignmentnumber" -> "TestData" ---
$ws.Name = "TestData"

# --- Drop the old A1:A11 contents ---
$ws.Cells.Clear()

# --- Write the new 4-column test-data table ---
$ws.Range("A1").Value = "TestCaseName"
$ws.Range("B1").Value = "ConsignmentNO"
$ws.Range("C1").Value = "ReceptacleID"
$ws.Range("D1").Value = "Status"

$ws.Range("A2").Value = "CommonTestData"
$ws.Range("B2").Value = "PREC01048239"
$ws.Range("C2").Value = "CARDITRECEPTACLEID00000111201"

$ws.Range("A3").Value = "AllOntime"
$ws.Range("B3").Value = "PREC01049379"
$ws.Range("D3").Value = "PRECON"

# --- Column widths the author had set on B, C, F, G ---
$ws.Columns.Item(2).ColumnWidth = 15.28515625
$ws.Columns.Item(3).ColumnWidth = 30.5703125
$ws.Columns.Item(6).ColumnWidth = 13.140625
$ws.Columns.Item(7).ColumnWidth = 30.85546875

# --- Selection / active sheet: TestData becomes the active tab,
#     selection parked one row below the table (D4) ---
[void]$ws.Range("D4").Select()
[void]$ws.Activate()
